$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# Overview sheet: status column now reflects "handed back" instead of
# "ready for handoff" for both rows.
# ---------------------------------------------------------------------------
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet: update status, fill in "Latest Target File" / "Latest Handback
# File" columns (F/G) that were generated for the handback report, and set
# the real handback datetime that replaces the zero-date placeholder.
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/94c238b74b5c0e66d59b2015c120b7cb13ac1f2d/e2e/34e5b160-d462-40e3-a381-0134ac95bb12.md", "", "", "34e5b160-d462-40e3-a381-0134ac95bb12.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e167bb772e9225f466fe5c8826781aa63ff9522/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/34e5b160-d462-40e3-a381-0134ac95bb12.f40c579fb2409bd58a9c038ccfd61cfd18672079.zh-cn.xlf", "", "", "34e5b160-d462-40e3-a381-0134ac95bb12.f40c579fb2409bd58a9c038ccfd61cfd18672079.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/94c238b74b5c0e66d59b2015c120b7cb13ac1f2d/e2e/5a765ba2-0408-43b4-afa6-c5f520b8e3dd.md", "", "", "5a765ba2-0408-43b4-afa6-c5f520b8e3dd.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e167bb772e9225f466fe5c8826781aa63ff9522/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/5a765ba2-0408-43b4-afa6-c5f520b8e3dd.84925c9c6e4f9e122f52f6f7d1362138e1297b8b.zh-cn.xlf", "", "", "5a765ba2-0408-43b4-afa6-c5f520b8e3dd.84925c9c6e4f9e122f52f6f7d1362138e1297b8b.zh-cn.xlf") | Out-Null

$zhCnNewCells = $wsZhCn.Range("F2:G3")
$zhCnNewCells.Font.Color = 15570276
$zhCnNewCells.Font.Underline = 2

$wsZhCn.Range("H2").Value = "2016-03-21 03:22:48"
$wsZhCn.Range("H3").Value = "2016-03-21 03:22:48"

# ---------------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but with its own handback
# datetime.
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/94c238b74b5c0e66d59b2015c120b7cb13ac1f2d/e2e/34e5b160-d462-40e3-a381-0134ac95bb12.md", "", "", "34e5b160-d462-40e3-a381-0134ac95bb12.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/438bbc7416af8faa4373410fd46b4efe644868ad/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/34e5b160-d462-40e3-a381-0134ac95bb12.f40c579fb2409bd58a9c038ccfd61cfd18672079.de-de.xlf", "", "", "34e5b160-d462-40e3-a381-0134ac95bb12.f40c579fb2409bd58a9c038ccfd61cfd18672079.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/94c238b74b5c0e66d59b2015c120b7cb13ac1f2d/e2e/5a765ba2-0408-43b4-afa6-c5f520b8e3dd.md", "", "", "5a765ba2-0408-43b4-afa6-c5f520b8e3dd.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/438bbc7416af8faa4373410fd46b4efe644868ad/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/5a765ba2-0408-43b4-afa6-c5f520b8e3dd.84925c9c6e4f9e122f52f6f7d1362138e1297b8b.de-de.xlf", "", "", "5a765ba2-0408-43b4-afa6-c5f520b8e3dd.84925c9c6e4f9e122f52f6f7d1362138e1297b8b.de-de.xlf") | Out-Null

$deDeNewCells = $wsDeDe.Range("F2:G3")
$deDeNewCells.Font.Color = 15570276
$deDeNewCells.Font.Underline = 2

$wsDeDe.Range("H2").Value = "2016-03-21 03:23:04"
$wsDeDe.Range("H3").Value = "2016-03-21 03:23:04"
